# Auto-generated: update market-board derived columns (H:N) for leve profit tables
# across multiple worksheets, per scheduled data-refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Cells.Item(28, 8).Value2 = 249.58333
$ws.Cells.Item(28, 9).Value2 = 217.72728
$ws.Cells.Item(28, 10).Value2 = 600
$ws.Cells.Item(28, 11).Value2 = 217.72728
$ws.Cells.Item(28, 12).Value2 = 600
$ws.Cells.Item(28, 13).Value2 = 267.27272
$ws.Cells.Item(28, 14).Value2 = -1570
# Row 113
$ws.Cells.Item(113, 8).Value2 = 2660.9
$ws.Cells.Item(113, 9).Value2 = 2708.8572
$ws.Cells.Item(113, 10).Value2 = 2635.077
$ws.Cells.Item(113, 11).Value2 = 2708.8572
$ws.Cells.Item(113, 12).Value2 = 2635.077
$ws.Cells.Item(113, 13).Value2 = 545.1428000000001
$ws.Cells.Item(113, 14).Value2 = -9143.077000000001
# Row 132
$ws.Cells.Item(132, 8).Value2 = 3341.5518
$ws.Cells.Item(132, 9).Value2 = 1793.8334
$ws.Cells.Item(132, 10).Value2 = 10770.6
$ws.Cells.Item(132, 11).Value2 = 5381.5002
$ws.Cells.Item(132, 12).Value2 = 32311.8
$ws.Cells.Item(132, 13).Value2 = -2851.5002
$ws.Cells.Item(132, 14).Value2 = -37371.8
# Row 137
$ws.Cells.Item(137, 8).Value2 = 2828.95
$ws.Cells.Item(137, 9).Value2 = 2851.75
$ws.Cells.Item(137, 10).Value2 = 2775.75
$ws.Cells.Item(137, 11).Value2 = 8555.25
$ws.Cells.Item(137, 12).Value2 = 8327.25
$ws.Cells.Item(137, 13).Value2 = -6005.25
$ws.Cells.Item(137, 14).Value2 = -13427.25
# Row 141
$ws.Cells.Item(141, 8).Value2 = 2226.6667
$ws.Cells.Item(141, 9).Value2 = 2226.6667
$ws.Cells.Item(141, 11).Value2 = 6680.000100000001
$ws.Cells.Item(141, 13).Value2 = -1500.000100000001

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value2 = 2925.5557
$ws.Cells.Item(2, 9).Value2 = 5550
$ws.Cells.Item(2, 10).Value2 = 826
$ws.Cells.Item(2, 11).Value2 = 5550
$ws.Cells.Item(2, 12).Value2 = 826
$ws.Cells.Item(2, 13).Value2 = -5437
$ws.Cells.Item(2, 14).Value2 = -1052
# Row 61
$ws.Cells.Item(61, 8).Value2 = 2422.5874
$ws.Cells.Item(61, 9).Value2 = 1588.5333
$ws.Cells.Item(61, 10).Value2 = 4507.722
$ws.Cells.Item(61, 11).Value2 = 1588.5333
$ws.Cells.Item(61, 12).Value2 = 4507.722
$ws.Cells.Item(61, 13).Value2 = -1376.5333
$ws.Cells.Item(61, 14).Value2 = -4931.722
# Row 81
$ws.Cells.Item(81, 8).Value2 = 18000
$ws.Cells.Item(81, 10).Value2 = 18000
$ws.Cells.Item(81, 12).Value2 = 18000
$ws.Cells.Item(81, 14).Value2 = -19996
# Row 84
$ws.Cells.Item(84, 8).Value2 = 18000
$ws.Cells.Item(84, 10).Value2 = 18000
$ws.Cells.Item(84, 12).Value2 = 54000
$ws.Cells.Item(84, 14).Value2 = -63984
# Row 116
$ws.Cells.Item(116, 8).Value2 = 2925.5557
$ws.Cells.Item(116, 9).Value2 = 5550
$ws.Cells.Item(116, 10).Value2 = 826
$ws.Cells.Item(116, 11).Value2 = 5550
$ws.Cells.Item(116, 12).Value2 = 826
$ws.Cells.Item(116, 13).Value2 = -3256
$ws.Cells.Item(116, 14).Value2 = -5414
# Row 122
$ws.Cells.Item(122, 8).Value2 = 3241.1052
$ws.Cells.Item(122, 9).Value2 = 3760.3333
$ws.Cells.Item(122, 10).Value2 = 2351
$ws.Cells.Item(122, 11).Value2 = 11280.9999
$ws.Cells.Item(122, 12).Value2 = 7053
$ws.Cells.Item(122, 13).Value2 = -8830.999899999999
$ws.Cells.Item(122, 14).Value2 = -11953
# Row 136
$ws.Cells.Item(136, 8).Value2 = 2422.5874
$ws.Cells.Item(136, 9).Value2 = 1588.5333
$ws.Cells.Item(136, 10).Value2 = 4507.722
$ws.Cells.Item(136, 11).Value2 = 4765.5999
$ws.Cells.Item(136, 12).Value2 = 13523.166
$ws.Cells.Item(136, 13).Value2 = -2215.5999
$ws.Cells.Item(136, 14).Value2 = -18623.166
# Row 140
$ws.Cells.Item(140, 8).Value2 = 42381.168
$ws.Cells.Item(140, 10).Value2 = 42381.168
$ws.Cells.Item(140, 12).Value2 = 42381.168
$ws.Cells.Item(140, 14).Value2 = -52741.168
# Row 141
$ws.Cells.Item(141, 8).Value2 = 30309.666
$ws.Cells.Item(141, 10).Value2 = 30309.666
$ws.Cells.Item(141, 12).Value2 = 30309.666
$ws.Cells.Item(141, 14).Value2 = -40669.666

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value2 = 2925.5557
$ws.Cells.Item(3, 9).Value2 = 5550
$ws.Cells.Item(3, 10).Value2 = 826
$ws.Cells.Item(3, 11).Value2 = 5550
$ws.Cells.Item(3, 12).Value2 = 826
$ws.Cells.Item(3, 13).Value2 = -5436
$ws.Cells.Item(3, 14).Value2 = -1054
# Row 86
$ws.Cells.Item(86, 8).Value2 = 7348.4443
$ws.Cells.Item(86, 9).Value2 = 5959.75
$ws.Cells.Item(86, 10).Value2 = 10125.833
$ws.Cells.Item(86, 11).Value2 = 5959.75
$ws.Cells.Item(86, 12).Value2 = 10125.833
$ws.Cells.Item(86, 13).Value2 = -4836.75
$ws.Cells.Item(86, 14).Value2 = -12371.833
# Row 89
$ws.Cells.Item(89, 8).Value2 = 7348.4443
$ws.Cells.Item(89, 9).Value2 = 5959.75
$ws.Cells.Item(89, 10).Value2 = 10125.833
$ws.Cells.Item(89, 11).Value2 = 29798.75
$ws.Cells.Item(89, 12).Value2 = 50629.165
$ws.Cells.Item(89, 13).Value2 = -24182.75
$ws.Cells.Item(89, 14).Value2 = -61861.165
# Row 132
$ws.Cells.Item(132, 8).Value2 = 33794.285
$ws.Cells.Item(132, 10).Value2 = 33794.285
$ws.Cells.Item(132, 12).Value2 = 33794.285
$ws.Cells.Item(132, 14).Value2 = -43914.285

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value2 = 2908.76
$ws.Cells.Item(31, 9).Value2 = 2116.1614
$ws.Cells.Item(31, 10).Value2 = 4201.9473
$ws.Cells.Item(31, 11).Value2 = 2116.1614
$ws.Cells.Item(31, 12).Value2 = 4201.9473
$ws.Cells.Item(31, 13).Value2 = -1821.1614
$ws.Cells.Item(31, 14).Value2 = -4791.9473
# Row 34
$ws.Cells.Item(34, 8).Value2 = 2908.76
$ws.Cells.Item(34, 9).Value2 = 2116.1614
$ws.Cells.Item(34, 10).Value2 = 4201.9473
$ws.Cells.Item(34, 11).Value2 = 2116.1614
$ws.Cells.Item(34, 12).Value2 = 4201.9473
$ws.Cells.Item(34, 13).Value2 = -1914.1614
$ws.Cells.Item(34, 14).Value2 = -4605.9473
# Row 58
$ws.Cells.Item(58, 8).Value2 = 1749.463
$ws.Cells.Item(58, 9).Value2 = 1307.1936
$ws.Cells.Item(58, 10).Value2 = 2345.5652
$ws.Cells.Item(58, 11).Value2 = 1307.1936
$ws.Cells.Item(58, 12).Value2 = 2345.5652
$ws.Cells.Item(58, 13).Value2 = -1104.1936
$ws.Cells.Item(58, 14).Value2 = -2751.5652
# Row 122
$ws.Cells.Item(122, 8).Value2 = 1108.4
$ws.Cells.Item(122, 9).Value2 = 743.4286
$ws.Cells.Item(122, 10).Value2 = 1960
$ws.Cells.Item(122, 11).Value2 = 2230.2858
$ws.Cells.Item(122, 12).Value2 = 5880
$ws.Cells.Item(122, 13).Value2 = 219.7142000000003
$ws.Cells.Item(122, 14).Value2 = -10780
# Row 132
$ws.Cells.Item(132, 8).Value2 = 2201.9038
$ws.Cells.Item(132, 9).Value2 = 1390.5454
$ws.Cells.Item(132, 10).Value2 = 3611.1052
$ws.Cells.Item(132, 11).Value2 = 4171.6362
$ws.Cells.Item(132, 12).Value2 = 10833.3156
$ws.Cells.Item(132, 13).Value2 = -1641.6362
$ws.Cells.Item(132, 14).Value2 = -15893.3156
# Row 134
$ws.Cells.Item(134, 8).Value2 = 1510.2885
$ws.Cells.Item(134, 9).Value2 = 1074.0605
$ws.Cells.Item(134, 10).Value2 = 2267.9473
$ws.Cells.Item(134, 11).Value2 = 3222.1815
$ws.Cells.Item(134, 12).Value2 = 6803.841899999999
$ws.Cells.Item(134, 13).Value2 = -687.1815000000001
$ws.Cells.Item(134, 14).Value2 = -11873.8419
# Row 136
$ws.Cells.Item(136, 8).Value2 = 1749.463
$ws.Cells.Item(136, 9).Value2 = 1307.1936
$ws.Cells.Item(136, 10).Value2 = 2345.5652
$ws.Cells.Item(136, 11).Value2 = 3921.5808
$ws.Cells.Item(136, 12).Value2 = 7036.6956
$ws.Cells.Item(136, 13).Value2 = -1371.5808
$ws.Cells.Item(136, 14).Value2 = -12136.6956

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value2 = 6549.674
$ws.Cells.Item(70, 9).Value2 = 3802.2632
$ws.Cells.Item(70, 10).Value2 = 19599.875
$ws.Cells.Item(70, 11).Value2 = 3802.2632
$ws.Cells.Item(70, 12).Value2 = 19599.875
$ws.Cells.Item(70, 13).Value2 = -3532.2632
$ws.Cells.Item(70, 14).Value2 = -20139.875
# Row 73
$ws.Cells.Item(73, 8).Value2 = 6549.674
$ws.Cells.Item(73, 9).Value2 = 3802.2632
$ws.Cells.Item(73, 10).Value2 = 19599.875
$ws.Cells.Item(73, 11).Value2 = 3802.2632
$ws.Cells.Item(73, 12).Value2 = 19599.875
$ws.Cells.Item(73, 13).Value2 = -2866.2632
$ws.Cells.Item(73, 14).Value2 = -21471.875
# Row 122
$ws.Cells.Item(122, 8).Value2 = 1953.3334
$ws.Cells.Item(122, 9).Value2 = 1776
$ws.Cells.Item(122, 10).Value2 = 2175
$ws.Cells.Item(122, 11).Value2 = 5328
$ws.Cells.Item(122, 12).Value2 = 6525
$ws.Cells.Item(122, 13).Value2 = -2878
$ws.Cells.Item(122, 14).Value2 = -11425
# Row 138
$ws.Cells.Item(138, 8).Value2 = 30679.8
$ws.Cells.Item(138, 10).Value2 = 30679.8
$ws.Cells.Item(138, 12).Value2 = 30679.8
$ws.Cells.Item(138, 14).Value2 = -40959.8
# Row 139
$ws.Cells.Item(139, 8).Value2 = 33715.6
$ws.Cells.Item(139, 10).Value2 = 33715.6
$ws.Cells.Item(139, 12).Value2 = 33715.6
$ws.Cells.Item(139, 14).Value2 = -43995.6
# Row 140
$ws.Cells.Item(140, 8).Value2 = 30118
$ws.Cells.Item(140, 9).Value2 = 10000
$ws.Cells.Item(140, 10).Value2 = 35147.5
$ws.Cells.Item(140, 11).Value2 = 10000
$ws.Cells.Item(140, 12).Value2 = 35147.5
$ws.Cells.Item(140, 14).Value2 = -45507.5
$ws.Cells.Item(140, 13).Value2 = -4820
# Row 141
$ws.Cells.Item(141, 8).Value2 = 46606.57
$ws.Cells.Item(141, 10).Value2 = 46606.57
$ws.Cells.Item(141, 12).Value2 = 46606.57
$ws.Cells.Item(141, 14).Value2 = -56966.57

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Cells.Item(46, 8).Value2 = 587
$ws.Cells.Item(46, 9).Value2 = 590.5
$ws.Cells.Item(46, 10).Value2 = 580
$ws.Cells.Item(46, 11).Value2 = 590.5
$ws.Cells.Item(46, 12).Value2 = 580
$ws.Cells.Item(46, 13).Value2 = -402.5
$ws.Cells.Item(46, 14).Value2 = -956
# Row 68
$ws.Cells.Item(68, 8).Value2 = 71432180
$ws.Cells.Item(68, 9).Value2 = 142858620
$ws.Cells.Item(68, 10).Value2 = 5717.143
$ws.Cells.Item(68, 11).Value2 = 142858620
$ws.Cells.Item(68, 12).Value2 = 5717.143
$ws.Cells.Item(68, 13).Value2 = -142857871
$ws.Cells.Item(68, 14).Value2 = -7215.143
# Row 71
$ws.Cells.Item(71, 8).Value2 = 71432180
$ws.Cells.Item(71, 9).Value2 = 142858620
$ws.Cells.Item(71, 10).Value2 = 5717.143
$ws.Cells.Item(71, 11).Value2 = 714293100
$ws.Cells.Item(71, 12).Value2 = 28585.715
$ws.Cells.Item(71, 13).Value2 = -714289356
$ws.Cells.Item(71, 14).Value2 = -36073.715
# Row 122
$ws.Cells.Item(122, 8).Value2 = 4257.143
$ws.Cells.Item(122, 9).Value2 = 4600
$ws.Cells.Item(122, 10).Value2 = 4000
$ws.Cells.Item(122, 11).Value2 = 13800
$ws.Cells.Item(122, 12).Value2 = 12000
$ws.Cells.Item(122, 13).Value2 = -11350
$ws.Cells.Item(122, 14).Value2 = -16900
# Row 132
$ws.Cells.Item(132, 8).Value2 = 6678.447
$ws.Cells.Item(132, 9).Value2 = 2347.5
$ws.Cells.Item(132, 10).Value2 = 12040.571
$ws.Cells.Item(132, 11).Value2 = 7042.5
$ws.Cells.Item(132, 12).Value2 = 36121.713
$ws.Cells.Item(132, 13).Value2 = -4512.5
$ws.Cells.Item(132, 14).Value2 = -41181.713

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Cells.Item(100, 8).Value2 = 10420
$ws.Cells.Item(100, 9).Value2 = 2504
$ws.Cells.Item(100, 10).Value2 = 50000
$ws.Cells.Item(100, 11).Value2 = 5008
$ws.Cells.Item(100, 12).Value2 = 100000
$ws.Cells.Item(100, 13).Value2 = -4467
$ws.Cells.Item(100, 14).Value2 = -101082
# Row 107
$ws.Cells.Item(107, 8).Value2 = 556.1
$ws.Cells.Item(107, 9).Value2 = 520.6667
$ws.Cells.Item(107, 10).Value2 = 609.25
$ws.Cells.Item(107, 11).Value2 = 1562.0001
$ws.Cells.Item(107, 12).Value2 = 1827.75
$ws.Cells.Item(107, 13).Value2 = 357.9999
$ws.Cells.Item(107, 14).Value2 = -5667.75
# Row 122
$ws.Cells.Item(122, 8).Value2 = 37402
$ws.Cells.Item(122, 9).Value2 = 38602.074
$ws.Cells.Item(122, 10).Value2 = 5000
$ws.Cells.Item(122, 11).Value2 = 115806.222
$ws.Cells.Item(122, 12).Value2 = 15000
$ws.Cells.Item(122, 13).Value2 = -113356.222
$ws.Cells.Item(122, 14).Value2 = -19900
# Row 136
$ws.Cells.Item(136, 8).Value2 = 1792.8254
$ws.Cells.Item(136, 9).Value2 = 1045.5883
$ws.Cells.Item(136, 10).Value2 = 2668.8965
$ws.Cells.Item(136, 11).Value2 = 3136.7649
$ws.Cells.Item(136, 12).Value2 = 8006.689499999999
$ws.Cells.Item(136, 13).Value2 = -586.7648999999997
$ws.Cells.Item(136, 14).Value2 = -13106.6895

